# The document holds a single 20-row x 5-column table. Only every 4th row
# (1, 5, 9, 13, 17) actually contains the "NN÷N=" expressions; the rows in
# between are blank spacer rows. Update each populated cell in place using
# Table.Cell(row, column) so duplicate expressions (e.g. "21÷2=", "48÷7=")
# are each mapped to their own distinct replacement instead of relying on a
# global Find/Replace that could mis-target repeated text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New values per content row, left-to-right (columns 1-5).
$rowValues = @{
    1  = @("86÷9=", "88÷2=", "27÷8=", "85÷3=", "53÷3=")
    5  = @("70÷9=", "43÷3=", "93÷3=", "53÷3=", "37÷9=")
    9  = @("76÷7=", "68÷4=", "53÷3=", "69÷6=", "66÷3=")
    13 = @("57÷8=", "13÷6=", "81÷5=", "72÷8=", "53÷8=")
    17 = @("49÷5=", "53÷9=", "41÷5=", "90÷8=", "55÷9=")
}

foreach ($row in $rowValues.Keys) {
    $values = $rowValues[$row]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
